$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.426422666666666
$ws.Range("H2").Value = 25.279268
$ws.Range("I2").Value = 0.1151758588783328
$ws.Range("J2").Value = 0.1151758588783328
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 948.0853966523518
$ws.Range("R2").Value = 8532.768569871167
$ws.Range("S2").Value = 0.03772369852637737
$ws.Range("T2").Value = 0.03772369852637737

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.426422666666666
$ws.Range("H3").Value = 25.279268
$ws.Range("I3").Value = 0.1151758588783328
$ws.Range("J3").Value = 0.1151758588783328
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 895.8506260969625
$ws.Range("R3").Value = 8062.655634872664
$ws.Range("S3").Value = 0.03564531113217879
$ws.Range("T3").Value = 0.0356453111321788

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.426422666666666
$ws.Range("H4").Value = 25.279268
$ws.Range("I4").Value = 0.1151758588783328
$ws.Range("J4").Value = 0.1151758588783328
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 1050.704590844987
$ws.Range("R4").Value = 9456.341317604887
$ws.Range("S4").Value = 0.04180684921977662
$ws.Range("T4").Value = 0.04180684921977663

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.78712033333333
$ws.Range("H5").Value = 35.361361
$ws.Range("I5").Value = 0.1611112760180311
$ws.Range("J5").Value = 0.1611112760180311
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 1326.208890615504
$ws.Range("R5").Value = 11935.88001553954
$ws.Range("S5").Value = 0.0527689853142266
$ws.Range("T5").Value = 0.0527689853142266

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.78712033333333
$ws.Range("H6").Value = 35.361361
$ws.Range("I6").Value = 0.1611112760180311
$ws.Range("J6").Value = 0.1611112760180311
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 1253.141403916075
$ws.Range("R6").Value = 11278.27263524468
$ws.Range("S6").Value = 0.04986167775515863
$ws.Range("T6").Value = 0.04986167775515864

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.78712033333333
$ws.Range("H7").Value = 35.361361
$ws.Range("I7").Value = 0.1611112760180311
$ws.Range("J7").Value = 0.1611112760180311
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 1469.755545976525
$ws.Range("R7").Value = 13227.79991378873
$ws.Range("S7").Value = 0.05848061294864589
$ws.Range("T7").Value = 0.05848061294864591

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 52.94781866666667
$ws.Range("H8").Value = 158.843456
$ws.Range("I8").Value = 0.7237128651036362
$ws.Range("J8").Value = 0.7237128651036362
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 5957.338677187584
$ws.Range("R8").Value = 53616.04809468825
$ws.Range("S8").Value = 0.2370386138962524
$ws.Range("T8").Value = 0.2370386138962524

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 52.94781866666667
$ws.Range("H9").Value = 158.843456
$ws.Range("I9").Value = 0.7237128651036362
$ws.Range("J9").Value = 0.7237128651036362
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 5629.119067411499
$ws.Range("R9").Value = 50662.07160670349
$ws.Range("S9").Value = 0.2239789700568289
$ws.Range("T9").Value = 0.223978970056829

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 52.94781866666667
$ws.Range("H10").Value = 158.843456
$ws.Range("I10").Value = 0.7237128651036362
$ws.Range("J10").Value = 0.7237128651036362
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 6602.151155835834
$ws.Range("R10").Value = 59419.3604025225
$ws.Range("S10").Value = 0.2626952811505548
$ws.Range("T10").Value = 0.2626952811505548
